# Updates crypto price/volume figures (and swaps the ShibaInu/Avalanche
# row order) to match the latest scrape.
# Note: some Price values look like plain numbers (e.g. "188.07"); a
# leading apostrophe is used so Excel stores them as text, matching the
# original inline-string cell type instead of silently coercing to a
# float (which would also mangle very small values like 0.0000187 into
# scientific notation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '74.913.61'
$ws.Range('E2').Value = '  +1.78%  '
$ws.Range('D3').Value = '2.820.43'
$ws.Range('E3').Value = '  +8.34%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''188.07'
$ws.Range('E5').Value = '  +1.88%  '
$ws.Range('D6').Value = '''596.23'
$ws.Range('E6').Value = '  +1.78%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '''0.556'
$ws.Range('E8').Value = '  +3.62%  '
$ws.Range('E9').Value = '  -4.53%  '
$ws.Range('D10').Value = '2.817.20'
$ws.Range('E10').Value = '  +8.50%  '
$ws.Range('E11').Value = '  -1.02%  '
$ws.Range('E12').Value = '  +2.09%  '
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('D14').Value = '3.335.39'
$ws.Range('E14').Value = '  +7.16%  '
$ws.Range('D15').Value = '74.818.91'
$ws.Range('E15').Value = '  +1.61%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '''0.0000187'
$ws.Range('E16').Value = '  -0.76%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').Value = '''26.96'
$ws.Range('E17').Value = '  +3.18%  '
$ws.Range('D18').Value = '2.817.92'
$ws.Range('E18').Value = '  +7.80%  '
$ws.Range('D19').Value = '''8.95'
$ws.Range('E19').Value = '  +0.71%  '
$ws.Range('D20').Value = '''12.30'
$ws.Range('E20').Value = '  +4.14%  '
$ws.Range('D21').Value = '''374.22'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('E22').Value = '  -0.50%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').Value = '''6.16'
$ws.Range('E24').Value = '  -0.61%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').Value = '''70.64'
$ws.Range('E26').Value = '  +1.06%  '
$ws.Range('D27').Value = '2.963.52'
$ws.Range('E27').Value = '  +8.43%  '
$ws.Range('D28').Value = '''4.15'
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('E29').Value = '  +2.75%  '
$ws.Range('D30').Value = '''0.0000102'
$ws.Range('E30').Value = '  +8.71%  '
$ws.Range('D31').Value = '''0.998'
$ws.Range('E31').Value = '  +0.25%  '
$ws.Range('D32').Value = '''511.81'
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('D33').Value = '''1.38'
$ws.Range('E33').Value = '  -0.10%  '
$ws.Range('D34').Value = '''7.84'
$ws.Range('E34').Value = '  -1.37%  '
$ws.Range('E35').Value = '  +2.76%  '
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('D37').Value = '''162.85'
$ws.Range('E37').Value = '  +1.92%  '
$ws.Range('E38').Value = '  +4.59%  '
$ws.Range('E39').Value = '  -2.22%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').Value = '''182.30'
$ws.Range('E41').Value = '  +15.89%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').Value = '''5.04'
$ws.Range('E43').Value = '  +2.88%  '
$ws.Range('E44').Value = '  +4.20%  '
$ws.Range('E45').Value = '  +0.35%  '
$ws.Range('E46').Value = '  +3.09%  '
$ws.Range('D47').Value = '''39.80'
$ws.Range('E47').Value = '  +2.88%  '
$ws.Range('E48').Value = '  -1.60%  '
$ws.Range('D49').Value = '''0.0865'
$ws.Range('E49').Value = '  -7.18%  '
$ws.Range('E50').Value = '  +7.80%  '
$ws.Range('E51').Value = '  +2.90%  '
